$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '53.892.27'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.252.41'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.69%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '493.84'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.93'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.523'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0945'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.36%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.335'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.73'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.651.15'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.89%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.62'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '53.841.41'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.17%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.63%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.252.42'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.19'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.12'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '300.23'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.29'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.00%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '60.63'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.998'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.148'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.25'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.22'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.50%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.91'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0682'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.08'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.938'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +7.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.19'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.69'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.370'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.39'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.35'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '124.79'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.77'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0489'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0888'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.540'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '238.79'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.369'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0204'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.75'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.08'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.86%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.03%  '
